# Refresh crypto price/volume table values (GitHub Actions scrape update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D holds price text that sometimes looks numeric (e.g. "0.9988", "1.000").
# Force those specific cells to Text format first so Excel keeps the exact string
# (preserves values like "1.000" / "6.840" / "4.750" instead of collapsing to 1 / 6.84 / 4.75).
$ws.Range("D4,D5,D6,D7,D8,D9,D10,D11,D12,D13,D14,D16,D17,D18,D19,D20,D21,D22,D24,D25,D26,D27,D29,D30,D31,D32,D34,D35,D36,D37,D38,D39,D40,D41,D42,D43,D44,D45,D46,D47,D48,D49,D50,D51").NumberFormat = "@"

# Updated Coin / Link / Price / Volume(1h) cells
$ws.Range("D2").Value = "27.138.44"
$ws.Range("E2").Value = "  -2.54%  "
$ws.Range("D3").Value = "1.715.62"
$ws.Range("E3").Value = "  -2.73%  "
$ws.Range("D4").Value = "0.9988"
$ws.Range("E4").Value = "  -0.38%  "
$ws.Range("D5").Value = "307.58"
$ws.Range("E5").Value = "  -6.20%  "
$ws.Range("D6").Value = "0.9998"
$ws.Range("E6").Value = "  -0.21%  "
$ws.Range("D7").Value = "0.4742"
$ws.Range("E7").Value = "  +6.45%  "
$ws.Range("D8").Value = "0.3434"
$ws.Range("E8").Value = "  -2.88%  "
$ws.Range("D9").Value = "41.94"
$ws.Range("E9").Value = "  -0.11%  "
$ws.Range("D10").Value = "0.07262"
$ws.Range("E10").Value = "  -1.85%  "
$ws.Range("D11").Value = "1.051"
$ws.Range("E11").Value = "  -4.19%  "
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  -0.15%  "
$ws.Range("D13").Value = "19.87"
$ws.Range("E13").Value = "  -4.93%  "
$ws.Range("D14").Value = "5.861"
$ws.Range("E14").Value = "  -2.61%  "
$ws.Range("D15").Value = "1.704.84"
$ws.Range("E15").Value = "  -3.28%  "
$ws.Range("D16").Value = "6.840"
$ws.Range("E16").Value = "  -5.36%  "
$ws.Range("D17").Value = "89.32"
$ws.Range("E17").Value = "  -3.91%  "
$ws.Range("D18").Value = "0.00001039"
$ws.Range("E18").Value = "  -2.05%  "
$ws.Range("D19").Value = "0.06360"
$ws.Range("E19").Value = "  -0.98%  "
$ws.Range("D20").Value = "1.000"
$ws.Range("D21").Value = "16.48"
$ws.Range("E21").Value = "  -3.53%  "
$ws.Range("D22").Value = "5.595"
$ws.Range("E22").Value = "  -2.82%  "
$ws.Range("D23").Value = "27.126.75"
$ws.Range("E23").Value = "  -2.73%  "
$ws.Range("D24").Value = "10.86"
$ws.Range("E24").Value = "  -3.40%  "
$ws.Range("D25").Value = "2.086"
$ws.Range("E25").Value = "  -1.07%  "
$ws.Range("D26").Value = "155.45"
$ws.Range("E26").Value = "  -3.97%  "
$ws.Range("D27").Value = "19.60"
$ws.Range("E27").Value = "  -3.74%  "
$ws.Range("D28").Value = "1.899.81"
$ws.Range("E28").Value = "  -3.40%  "
$ws.Range("D29").Value = "2.082"
$ws.Range("E29").Value = "  -3.38%  "
$ws.Range("D30").Value = "119.65"
$ws.Range("E30").Value = "  -4.15%  "
$ws.Range("D31").Value = "1.014"
$ws.Range("E31").Value = "  -7.81%  "
$ws.Range("D32").Value = "0.09166"
$ws.Range("E32").Value = "  -0.21%  "
$ws.Range("E33").Value = "  -2.83%  "
$ws.Range("D34").Value = "5.325"
$ws.Range("E34").Value = "  -5.38%  "
$ws.Range("D35").Value = "0.02196"
$ws.Range("E35").Value = "  -3.57%  "
$ws.Range("D36").Value = "0.05829"
$ws.Range("E36").Value = "  -5.84%  "
$ws.Range("D37").Value = "11.09"
$ws.Range("E37").Value = "  -6.30%  "
$ws.Range("B38").Value = "Algorand"
$ws.Range("C38").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D38").Value = "0.1995"
$ws.Range("E38").Value = "  -4.86%  "
$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").Value = "4.750"
$ws.Range("E39").Value = "  -3.98%  "
$ws.Range("B40").Value = "WEMIXTOKEN"
$ws.Range("C40").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D40").Value = "1.393"
$ws.Range("E40").Value = "  -0.01%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "0.5851"
$ws.Range("E41").Value = "  -7.31%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "1.102"
$ws.Range("E42").Value = "  -7.15%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "7.480"
$ws.Range("E43").Value = "  -4.84%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "12.65"
$ws.Range("E44").Value = "  -4.64%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "0.5645"
$ws.Range("E45").Value = "  -3.61%  "
$ws.Range("D46").Value = "3.553"
$ws.Range("E46").Value = "  -5.22%  "
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").Value = "117.64"
$ws.Range("E47").Value = "  -4.07%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "1.837"
$ws.Range("E48").Value = "  -5.93%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "0.06633"
$ws.Range("E49").Value = "  -3.80%  "
$ws.Range("B50").Value = "EOS"
$ws.Range("C50").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D50").Value = "1.085"
$ws.Range("E50").Value = "  -4.27%  "
$ws.Range("B51").Value = "PaxDollar"
$ws.Range("C51").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D51").Value = "0.9999"
$ws.Range("E51").Value = "  -0.17%  "
